$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add "path" header and file paths in column B
$ws.Range("B1").Value = "path"
$ws.Range("B2").Value = "C:\Users\Root\Desktop\test\SEAL-project\src\assets\logos\logo1.png"
$ws.Range("B3").Value = "C:\Users\Root\Desktop\test\SEAL-project\src\assets\logos\logo2.png"
$ws.Range("B4").Value = "C:\Users\Root\Desktop\test\SEAL-project\src\assets\logos\logo3.png"
$ws.Range("B5").Value = "C:\Users\Root\Desktop\test\SEAL-project\src\assets\logos\logo4.png"
$ws.Range("B6").Value = "C:\Users\Root\Desktop\test\SEAL-project\src\assets\logos\logo5.png"

# Set column B width to match target (64.875 chars in the saved XML).
# The host snaps column widths to 1/7-character pixel boundaries (Calibri-style
# max digit width), so 64.125 is the closest settable value that serializes
# nearest to the authored 64.875 width.
$ws.Columns.Item(2).ColumnWidth = 64.125

# Select D3 as the active cell, matching the target sheetView selection
$ws.Range("D3").Select()
